$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'35.281.35"
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').Value = "'1.897.08"
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'0.693"
$ws.Range('E5').Value = '  +9.30%  '
$ws.Range('D6').Value = "'245.08"
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = "'40.72"
$ws.Range('E8').Value = '  -4.15%  '
$ws.Range('D9').Value = "'0.347"
$ws.Range('E9').Value = '  +2.30%  '
$ws.Range('D10').Value = "'53.00"
$ws.Range('E10').Value = '  +11.01%  '
$ws.Range('D11').Value = "'0.0719"
$ws.Range('E11').Value = '  +1.84%  '
$ws.Range('D12').Value = "'0.0994"
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').Value = "'2.174.93"
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('D14').Value = "'12.43"
$ws.Range('E14').Value = '  -0.55%  '
$ws.Range('D15').Value = "'0.703"
$ws.Range('E15').Value = '  +1.70%  '
$ws.Range('D16').Value = "'1.893.81"
$ws.Range('E16').Value = '  -0.53%  '
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').Value = "'35.279.68"
$ws.Range('E18').Value = '  -0.77%  '
$ws.Range('D19').Value = "'72.11"
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').Value = "'0.0₃0816"
$ws.Range('E20').Value = '  +0.83%  '
$ws.Range('D21').Value = "'240.40"
$ws.Range('E21').Value = '  -1.65%  '
$ws.Range('D22').Value = "'12.57"
$ws.Range('E22').Value = '  +0.96%  '
$ws.Range('D23').Value = "'4.78"
$ws.Range('E23').Value = '  -2.70%  '
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('D25').Value = "'2.31"
$ws.Range('E25').Value = '  +1.03%  '
$ws.Range('E26').Value = '  +7.51%  '
$ws.Range('D27').Value = "'167.95"
$ws.Range('E27').Value = '  -1.84%  '
$ws.Range('D28').Value = "'8.57"
$ws.Range('E28').Value = '  +1.49%  '
$ws.Range('E29').Value = '  +3.90%  '
$ws.Range('E30').Value = '  +2.02%  '
$ws.Range('D31').Value = "'4.142.54"
$ws.Range('E31').Value = '  +21.33%  '
$ws.Range('E32').Value = '  +1.17%  '
$ws.Range('E33').Value = '  +0.19%  '
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('D35').Value = "'0.920"
$ws.Range('E35').Value = '  -2.92%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').Value = "'1.52"
$ws.Range('E36').Value = '  +14.14%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = "'4.09"
$ws.Range('E37').Value = '  -0.52%  '
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').Value = "'1.81"
$ws.Range('E38').Value = '  +2.42%  '
$ws.Range('D39').Value = "'2.01"
$ws.Range('E39').Value = '  -1.38%  '
$ws.Range('E40').Value = '  +10.83%  '
$ws.Range('E41').Value = '  +2.11%  '
$ws.Range('E42').Value = '  -1.18%  '
$ws.Range('E43').Value = '  +4.73%  '
$ws.Range('D44').Value = "'89.45"
$ws.Range('E44').Value = '  -2.01%  '
$ws.Range('D45').Value = "'1.350.17"
$ws.Range('E45').Value = '  -0.56%  '
$ws.Range('D46').Value = "'2.41"
$ws.Range('E46').Value = '  +2.59%  '
$ws.Range('E47').Value = '  -2.68%  '
$ws.Range('D48').Value = "'2.43"
$ws.Range('E48').Value = '  +0.29%  '
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('D50').Value = "'45.93"
$ws.Range('E50').Value = '  -2.13%  '
$ws.Range('D51').Value = "'6.50"
$ws.Range('E51').Value = '  -2.52%  '
